$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 used to only have A12 ("21 loka") and B12 (a time range) filled in.
# This entry is being fleshed out into a full diary row like the others:
# the old time text moves out of B12's slot conceptually and a freshly
# logged (slightly later-finishing) time range takes its place, while the
# rest of the row (content / quality / notes / meta / hours) gets filled in.

$ws.Range("B12").Value = "9.15-10.30, 14.00-17.00, 19.00-19.45"

$ws.Range("C12").Value = "Infinite grid yritelmää"
$ws.Range("D12").Value = "Aikalailla täysi nolla, kunnes tajusin että niin, klippaushan siinä tapahtuu. Kaikkia tunteja ei ole viitsitty kirjata kun eivät olleet tehokkaita."
$ws.Range("E12").Value = "Tuli pientä sörkkimistä joka puolelta debug mielessä, mutta periaatteisiinhan se kosahti. Near - plane oli liian lähellä, ja pieni ""hack"" joka nosti y tasoa hieman selitti enemmän kuin haluan myöntää."
$ws.Range("F12").Value = "Täytyy alkaa pohtia jos alkaa oikeasti mennä näin kauan aikaa ""perustan"" parantamiseen ja openGL kikkailuihin, onko mielekästä kurssia tässä kohtaa yrittää käydä? Priorisoidaan, ja onhan tässä 10 viikkoa. Toivottavasti pian päästään kirjassakin eteenpäin."
$ws.Range("G12").Value = 5

# Match the wrap-text formatting used by the other narrative columns.
$ws.Range("C12:F12").WrapText = $true

# The row grows a lot taller now that it holds full paragraphs of text.
$ws.Rows.Item(12).RowHeight = 116

# Leave the cursor where the author ended up after typing the new entry.
$ws.Range("F13").Select()
